$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-9, columns A:T
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric metrics
$data = @(
  @("ECs","Tnfsf14","Ltbr","ECs",1,0.3333333333333333,0.02601766666666666,0.078053,0.02802616153292364,0.02802616153292364,3,1,11.30399166666667,33.911975,0.2607347834842907,0.2607347834842906,0.2941034871861111,2.646931384675,0.007307395159182601,0.0073073951591826),
  @("ECs","Tnfsf14","Ltbr","FAPs",1,0.3333333333333333,0.02601766666666666,0.078053,0.02802616153292364,0.02802616153292364,3,1,17.704387,53.11316100000001,0.4083645536274804,0.4083645536274804,0.4606268395036666,4.145641555533,0.01144489094428402,0.01144489094428402),
  @("ECs","Tnfsf14","Ltbr","M2",1,0.3333333333333333,0.02601766666666666,0.078053,0.02802616153292364,0.02802616153292364,3,1,8.574736333333332,25.724209,0.1977825256098956,0.1977825256098956,0.2230946316752222,2.007851685077,0.005543085011132542,0.005543085011132542),
  @("ECs","Tnfsf14","Ltbr","sCs",1,0.3333333333333333,0.02601766666666666,0.078053,0.02802616153292364,0.02802616153292364,3,1,5.771252666666666,17.313758,0.1331181372783333,0.1331181372783333,0.1501545281304444,1.351390753174,0.003730790418324472,0.003730790418324472),
  @("M2","Tnfsf14","Ltbr","ECs",3,1,0.9023173333333333,2.706952,0.9719738384670763,0.9719738384670763,3,1,11.30399166666667,33.911975,0.2607347834842907,0.2607347834842906,10.19978761668889,91.79808855019999,0.2534273883251081,0.253427388325108),
  @("M2","Tnfsf14","Ltbr","FAPs",3,1,0.9023173333333333,2.706952,0.9719738384670763,0.9719738384670763,3,1,17.704387,53.11316100000001,0.4083645536274804,0.4083645536274804,15.97497526614133,143.774777395272,0.3969196626831963,0.3969196626831963),
  @("M2","Tnfsf14","Ltbr","M2",3,1,0.9023173333333333,2.706952,0.9719738384670763,0.9719738384670763,3,1,8.574736333333332,25.724209,0.1977825256098956,0.1977825256098956,7.737133222329777,69.634199000968,0.1922394405987631,0.1922394405987631),
  @("M2","Tnfsf14","Ltbr","sCs",3,1,0.9023173333333333,2.706952,0.9719738384670763,0.9719738384670763,3,1,5.771252666666666,17.313758,0.1331181372783333,0.1331181372783333,5.207501316179555,46.86751184561599,0.1293873468600088,0.1293873468600088)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  $rowData = $data[$i]
  for ($c = 0; $c -lt $rowData.Count; $c++) {
    $col = $c + 1
    $ws.Cells.Item($row, $col).Value = $rowData[$c]
  }
}
